$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Architecture")
Write-Host $ws.Name
